# Lecturer_Class_Mapping.xlsx
# - Update the lecturer e-mail for the "Microcontroller/1B", "Full Stack
#   Development/1A" and "Full Stack Development/1B" rows (C3, C6, C7) from
#   2401404e@gmail.com to junleiliew@gmail.com.
# - Add individual mailto hyperlinks on those same cells (in the order the
#   author clicked them: C6, then C3, then C7).
# - Move the selection to C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newEmail = "junleiliew@gmail.com"

# --- Update the displayed e-mail text for the three "risky" rows ---
$ws.Range("C3").Value = $newEmail
$ws.Range("C6").Value = $newEmail
$ws.Range("C7").Value = $newEmail

# --- Add the new per-cell hyperlinks (author order: C6, C3, C7) ---
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:" + $newEmail)
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:" + $newEmail)
$ws.Hyperlinks.Add($ws.Range("C7"), "mailto:" + $newEmail)

# Re-apply the (pre-existing) Hyperlink cell style so the visible formatting
# matches what it was before (Hyperlinks.Add otherwise stamps a fresh style
# record); C4 still carries the original, untouched Hyperlink style.
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C3").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C6").PasteSpecial(-4122) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# --- Move the active selection to C3 ---
$ws.Range("C3").Select() | Out-Null

Write-Host "Lecturer e-mail updated and hyperlinks added."
